$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.742.09"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "2.248.68"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'249.04"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "'70.47"
$ws.Range("E7").Value = "  +6.67%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.660"
$ws.Range("E9").Value = "  +16.23%  "
$ws.Range("D10").Value = "'38.94"
$ws.Range("E10").Value = "  +8.98%  "
$ws.Range("D11").Value = "'59.51"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "'0.0953"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  +7.84%  "
$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "2.577.86"
$ws.Range("D16").Value = "'14.74"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "'0.875"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "2.254.84"
$ws.Range("E18").Value = "  +4.31%  "
$ws.Range("D19").Value = "42.674.74"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'234.47"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  +6.07%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'11.41"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "'2.29"
$ws.Range("E30").Value = "  +15.02%  "
$ws.Range("D31").Value = "'167.58"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "'6.49"
$ws.Range("E33").Value = "  +14.17%  "
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").Value = "'0.0797"
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("D36").Value = "'31.29"
$ws.Range("E36").Value = "  +25.91%  "
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = "  +10.89%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  +7.91%  "
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").Value = "'5.79"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'12.31"
$ws.Range("E43").Value = "  +6.20%  "
$ws.Range("D44").Value = "'62.10"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").Value = "'8.96"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("E51").Value = "  +3.54%  "
